$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 58564
$ws.Range("B4").Value = 58256
$ws.Range("B5").Value = 58043
$ws.Range("B6").Value = 58564
